# Append a new daily data row (row 42) to each of the four log sheets,
# mirroring the existing "time / length / ID / actual length / checksum"
# table layout used by rows 2-41.

$wb = $excel.ActiveWorkbook

# Per-sheet values for the new row 42.
# A = timestamp (days, Excel serial date)
# B = total-length hex bytes
# C = ID hex bytes
# D = actual-length hex bytes
# E = checksum hex
# F..I = decimal counterparts of B..E
$rowsToAdd = @{
    "FE_LFT_#1" = @{
        A = "45828.49421296296"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x60"
        E = "0xf"
        F = "380"
        G = "7.598631275147109e+23"
        H = "352"
        I = "15"
    }
    "FE_LFT_#2" = @{
        A = "45828.49421296296"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x74"
        E = "0xe"
        F = "400"
        G = "5.68432987514711e+23"
        H = "372"
        I = "14"
    }
    "FE_PLT_#1" = @{
        A = "45828.49421296296"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6A"
        E = "0x3"
        F = "110"
        G = "5.68631262647114e+23"
        H = "106"
        I = "3"
    }
    "FE_PLT_#2" = @{
        A = "45828.49421296296"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6A"
        E = "0x3"
        F = "110"
        G = "9.85046333984776e+23"
        H = "106"
        I = "3"
    }
}

foreach ($sheetName in $rowsToAdd.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 42
    $prevRow = $newRow - 1
    $data = $rowsToAdd[$sheetName]

    # Match the date/time number format used by the preceding rows in column A.
    $ws.Range("A$newRow").NumberFormat = $ws.Range("A$prevRow").NumberFormat
    $ws.Range("A$newRow").Value = [double]$data.A

    $ws.Range("B$newRow").Value = $data.B
    $ws.Range("C$newRow").Value = $data.C
    $ws.Range("D$newRow").Value = $data.D
    $ws.Range("E$newRow").Value = $data.E

    $ws.Range("F$newRow").Value = [double]$data.F
    $ws.Range("G$newRow").Value = [double]$data.G
    $ws.Range("H$newRow").Value = [double]$data.H
    $ws.Range("I$newRow").Value = [double]$data.I
}
